$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 99999
$ws.Cells.Item(62, 9).Value = 99999
$ws.Cells.Item(62, 11).Value = 99999
$ws.Cells.Item(62, 13).Value = -99375
$ws.Cells.Item(65, 8).Value = 99999
$ws.Cells.Item(65, 9).Value = 99999
$ws.Cells.Item(65, 11).Value = 499995
$ws.Cells.Item(65, 13).Value = -496875
$ws.Cells.Item(125, 8).Value = 1327.4
$ws.Cells.Item(125, 9).Value = 1327.4
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 11).Value = 11946.6
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 13).Value = -9486.6
$ws.Cells.Item(125, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 50005412
$ws.Cells.Item(132, 9).Value = 52637064
$ws.Cells.Item(132, 10).Value = 3999
$ws.Cells.Item(132, 11).Value = 157911192
$ws.Cells.Item(132, 12).Value = 11997
$ws.Cells.Item(132, 13).Value = -157908662
$ws.Cells.Item(132, 14).Value = -17057
$ws.Cells.Item(137, 8).Value = 19609700
$ws.Cells.Item(137, 9).Value = 30304264
$ws.Cells.Item(137, 11).Value = 90912792
$ws.Cells.Item(137, 13).Value = -90910242

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2722
$ws.Cells.Item(2, 9).Value = 2679.5
$ws.Cells.Item(2, 11).Value = 2679.5
$ws.Cells.Item(2, 13).Value = -2566.5
$ws.Cells.Item(32, 8).Value = 9815.0625
$ws.Cells.Item(32, 9).Value = 8617.566000000001
$ws.Cells.Item(32, 10).Value = 27777.5
$ws.Cells.Item(32, 11).Value = 8617.566000000001
$ws.Cells.Item(32, 12).Value = 27777.5
$ws.Cells.Item(32, 13).Value = -8330.566000000001
$ws.Cells.Item(32, 14).Value = -28351.5
$ws.Cells.Item(61, 8).Value = 7831.364
$ws.Cells.Item(61, 9).Value = 8119.5
$ws.Cells.Item(61, 11).Value = 8119.5
$ws.Cells.Item(61, 13).Value = -7907.5
$ws.Cells.Item(74, 8).Value = 3452273.5
$ws.Cells.Item(74, 9).Value = 4060734.8
$ws.Cells.Item(74, 10).Value = 4326
$ws.Cells.Item(74, 11).Value = 4060734.8
$ws.Cells.Item(74, 12).Value = 4326
$ws.Cells.Item(74, 13).Value = -4059860.8
$ws.Cells.Item(74, 14).Value = -6074
$ws.Cells.Item(77, 8).Value = 3452273.5
$ws.Cells.Item(77, 9).Value = 4060734.8
$ws.Cells.Item(77, 10).Value = 4326
$ws.Cells.Item(77, 11).Value = 20303674
$ws.Cells.Item(77, 12).Value = 21630
$ws.Cells.Item(77, 13).Value = -20299306
$ws.Cells.Item(77, 14).Value = -30366
$ws.Cells.Item(102, 8).Value = 7399.2856
$ws.Cells.Item(102, 9).Value = 6299.25
$ws.Cells.Item(102, 10).Value = 8866
$ws.Cells.Item(102, 11).Value = 6299.25
$ws.Cells.Item(102, 12).Value = 8866
$ws.Cells.Item(102, 13).Value = -4677.25
$ws.Cells.Item(102, 14).Value = -12110
$ws.Cells.Item(116, 8).Value = 2722
$ws.Cells.Item(116, 9).Value = 2679.5
$ws.Cells.Item(116, 11).Value = 2679.5
$ws.Cells.Item(116, 13).Value = -385.5
$ws.Cells.Item(122, 8).Value = 3495.5789
$ws.Cells.Item(122, 9).Value = 3115.9412
$ws.Cells.Item(122, 11).Value = 9347.8236
$ws.Cells.Item(122, 13).Value = -6897.8236
$ws.Cells.Item(126, 8).Value = 0
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 13).ClearContents()
$ws.Cells.Item(131, 8).Value = 78750
$ws.Cells.Item(131, 10).Value = 78750
$ws.Cells.Item(131, 12).Value = 78750
$ws.Cells.Item(131, 14).Value = -88830
$ws.Cells.Item(132, 8).Value = 5305.56
$ws.Cells.Item(132, 9).Value = 5305.56
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 15916.68
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -13386.68
$ws.Cells.Item(132, 14).ClearContents()
$ws.Cells.Item(136, 8).Value = 7831.364
$ws.Cells.Item(136, 9).Value = 8119.5
$ws.Cells.Item(136, 11).Value = 24358.5
$ws.Cells.Item(136, 13).Value = -21808.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2722
$ws.Cells.Item(3, 9).Value = 2679.5
$ws.Cells.Item(3, 11).Value = 2679.5
$ws.Cells.Item(3, 13).Value = -2565.5
$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 13).ClearContents()
$ws.Cells.Item(108, 8).Value = 69999.25
$ws.Cells.Item(108, 10).Value = 69999.25
$ws.Cells.Item(108, 12).Value = 69999.25
$ws.Cells.Item(108, 14).Value = -77679.25
$ws.Cells.Item(134, 8).Value = 2006.7273
$ws.Cells.Item(134, 9).Value = 897.44446
$ws.Cells.Item(134, 10).Value = 6998.5
$ws.Cells.Item(134, 11).Value = 2692.33338
$ws.Cells.Item(134, 12).Value = 20995.5
$ws.Cells.Item(134, 13).Value = -157.33338
$ws.Cells.Item(134, 14).Value = -26065.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2502.4
$ws.Cells.Item(31, 9).Value = 2428
$ws.Cells.Item(31, 11).Value = 2428
$ws.Cells.Item(31, 13).Value = -2133
$ws.Cells.Item(34, 8).Value = 2502.4
$ws.Cells.Item(34, 9).Value = 2428
$ws.Cells.Item(34, 11).Value = 2428
$ws.Cells.Item(34, 13).Value = -2226
$ws.Cells.Item(58, 8).Value = 1944.8572
$ws.Cells.Item(58, 9).Value = 1832.1177
$ws.Cells.Item(58, 11).Value = 1832.1177
$ws.Cells.Item(58, 13).Value = -1629.1177
$ws.Cells.Item(122, 8).Value = 5084.1
$ws.Cells.Item(122, 9).Value = 6232.7144
$ws.Cells.Item(122, 10).Value = 2404
$ws.Cells.Item(122, 11).Value = 18698.1432
$ws.Cells.Item(122, 12).Value = 7212
$ws.Cells.Item(122, 13).Value = -16248.1432
$ws.Cells.Item(122, 14).Value = -12112
$ws.Cells.Item(132, 8).Value = 11773246
$ws.Cells.Item(132, 9).Value = 11773246
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 35319738
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -35317208
$ws.Cells.Item(132, 14).ClearContents()
$ws.Cells.Item(136, 8).Value = 1944.8572
$ws.Cells.Item(136, 9).Value = 1832.1177
$ws.Cells.Item(136, 11).Value = 5496.3531
$ws.Cells.Item(136, 13).Value = -2946.3531

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38, 8).Value = 56.57143
$ws.Cells.Item(38, 9).Value = 59
$ws.Cells.Item(38, 10).Value = 50.5
$ws.Cells.Item(38, 11).Value = 177
$ws.Cells.Item(38, 12).Value = 151.5
$ws.Cells.Item(38, 13).Value = 170
$ws.Cells.Item(38, 14).Value = -845.5
$ws.Cells.Item(113, 8).Value = 3545.3333
$ws.Cells.Item(113, 9).Value = 3950
$ws.Cells.Item(113, 10).Value = 2736
$ws.Cells.Item(113, 11).Value = 11850
$ws.Cells.Item(113, 12).Value = 8208
$ws.Cells.Item(113, 13).Value = -9680
$ws.Cells.Item(113, 14).Value = -12548

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 4457.6665
$ws.Cells.Item(122, 9).Value = 4528.9165
$ws.Cells.Item(122, 11).Value = 13586.7495
$ws.Cells.Item(122, 13).Value = -11136.7495
$ws.Cells.Item(126, 8).Value = 0
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 13).ClearContents()
$ws.Cells.Item(132, 8).Value = 47621524
$ws.Cells.Item(132, 9).Value = 2885.6667
$ws.Cells.Item(132, 10).Value = 333333340
$ws.Cells.Item(132, 11).Value = 8657.000100000001
$ws.Cells.Item(132, 12).Value = 1000000020
$ws.Cells.Item(132, 13).Value = -6127.000100000001
$ws.Cells.Item(132, 14).Value = -1000005080

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 125000950
$ws.Cells.Item(22, 10).Value = 500000000
$ws.Cells.Item(22, 12).Value = 500000000
$ws.Cells.Item(22, 14).Value = -500000590
$ws.Cells.Item(27, 8).Value = 125000950
$ws.Cells.Item(27, 10).Value = 500000000
$ws.Cells.Item(27, 12).Value = 500000000
$ws.Cells.Item(27, 14).Value = -500000214
$ws.Cells.Item(122, 8).Value = 3325.111
$ws.Cells.Item(122, 9).Value = 3115.75
$ws.Cells.Item(122, 11).Value = 9347.25
$ws.Cells.Item(122, 13).Value = -6897.25
$ws.Cells.Item(132, 8).Value = 5297.0435
$ws.Cells.Item(132, 9).Value = 3434.9375
$ws.Cells.Item(132, 10).Value = 9553.286
$ws.Cells.Item(132, 11).Value = 10304.8125
$ws.Cells.Item(132, 12).Value = 28659.858
$ws.Cells.Item(132, 13).Value = -7774.8125
$ws.Cells.Item(132, 14).Value = -33719.858
$ws.Cells.Item(136, 8).Value = 2247.1177
$ws.Cells.Item(136, 9).Value = 2206.3125
$ws.Cells.Item(136, 11).Value = 6618.9375
$ws.Cells.Item(136, 13).Value = -4068.9375

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(109, 8).Value = 60000
$ws.Cells.Item(109, 10).Value = 60000
$ws.Cells.Item(109, 12).Value = 60000
$ws.Cells.Item(109, 14).Value = -62774
$ws.Cells.Item(122, 8).Value = 2146.9614
$ws.Cells.Item(122, 9).Value = 2146.9614
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 6440.8842
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -3990.8842
$ws.Cells.Item(122, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 166667630
$ws.Cells.Item(132, 9).Value = 799.6667
$ws.Cells.Item(132, 11).Value = 2399.0001
$ws.Cells.Item(132, 13).Value = 130.9998999999998
$ws.Cells.Item(133, 8).Value = 121248.25
$ws.Cells.Item(133, 10).Value = 121248.25
$ws.Cells.Item(133, 12).Value = 121248.25
$ws.Cells.Item(133, 14).Value = -131368.25
$ws.Cells.Item(136, 8).Value = 5565.1035
$ws.Cells.Item(136, 9).Value = 5648.852
$ws.Cells.Item(136, 11).Value = 16946.556
$ws.Cells.Item(136, 13).Value = -14396.556
